# Refresh the crypto price/volume snapshot columns (D = Price, E = Volume(1h))
# with the latest scraped values, as produced by the symbol-list update job.
# Cells are forced to Text format first so Excel keeps the values exactly as
# strings (matching how the source data is stored) instead of re-interpreting
# them as numbers/percentages.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "245.53"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-0.58%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "28.23"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-3.76%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.280"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "1.71%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05710"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-0.42%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.639"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "1.44%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.210"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "3.63%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8506"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-1.04%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.8858"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "1.78%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1382"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "1.32%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07086"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "0.15%"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "5.40%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09222"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-1.72%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001526"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.98%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0005958"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-1.16%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.006062"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "1.19%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.492"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "0.04%"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-4.01%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.03320"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-2.20%"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-0.40%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.533"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "1.85%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04077"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-1.16%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.1379"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.001219"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-0.50%"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-17.03%"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-0.85%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03789"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "0.95%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1068"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-0.34%"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-34.50%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002389"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-1.60%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.009482"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "10.64%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005271"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "0.38%"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-0.01%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1050"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "84.22%"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-0.34%"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-0.01%"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.01%"
